$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting (style + row height template) from the last existing
# row (14) down into the newly added rows (15-19) so the new cells pick up
# the same style index (wrap text, font) as the rest of the sheet. ---
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B19").PasteSpecial(-4122)

# --- Column A grows wider to fit the new "Powerful Positive Motivation
# Quotes" keyword text. ---
$ws.Columns.Item(1).ColumnWidth = 27.08

# --- Row 1 (header) is unchanged. ---
$ws.Range("A1").Value = "keyword"
$ws.Range("B1").Value = "appID"

# --- Row 2: new "Powerful Positive Motivation Quotes" entry inserted. ---
$ws.Range("A2").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B2").Value = "com.sugar.powerfulquotes"

# --- Row 3 ---
$ws.Range("A3").Value = "earn passive income"
$ws.Range("B3").Value = "passive.income.nadi.myfirstdrawermenuproject2"

# --- Row 4 ---
$ws.Range("A4").Value = "bitcoin"
$ws.Range("B4").Value = "com.hamxa.shaynachim"

# --- Row 5 ---
$ws.Range("A5").Value = "passive income ideas"
$ws.Range("B5").Value = "passive.income.nadi.myfirstdrawermenuproject2"

# --- Row 6 ---
$ws.Range("A6").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B6").Value = "com.sugar.powerfulquotes"

# --- Row 7 ---
$ws.Range("A7").Value = "best bitcoin"
$ws.Range("B7").Value = "com.hamxa.shaynachim"

# --- Row 8 ---
$ws.Range("A8").Value = "bitcoin beginners  "
$ws.Range("B8").Value = "com.hamxa.shaynachim"

# --- Row 9 ---
$ws.Range("A9").Value = "earn passive income"
$ws.Range("B9").Value = "passive.income.nadi.myfirstdrawermenuproject2"

# --- Row 10 ---
$ws.Range("A10").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B10").Value = "com.sugar.powerfulquotes"

# --- Row 11 ---
$ws.Range("A11").Value = "Best bitcoin"
$ws.Range("B11").Value = "com.hamxa.shaynachim"

# --- Row 12 ---
$ws.Range("A12").Value = "bitcoin"
$ws.Range("B12").Value = "com.hamxa.shaynachim"

# --- Row 13 ---
$ws.Range("A13").Value = "bitcoin"
$ws.Range("B13").Value = "com.hamxa.shaynachim"

# --- Row 14 ---
$ws.Range("A14").Value = "blockchain technology"
$ws.Range("B14").Value = "block.chain.technology"

# --- Row 15 (new): keyword column left blank. ---
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = "com.sugar.powerfulquotes"

# --- Row 16 (new) ---
$ws.Range("A16").Value = "bitcoin"
$ws.Range("B16").Value = "com.hamxa.shaynachim"

# --- Row 17 (new) ---
$ws.Range("A17").Value = "bitcoin"
$ws.Range("B17").Value = "com.hamxa.shaynachim"

# --- Row 18 (new) ---
$ws.Range("A18").Value = "affiliate marketing"
$ws.Range("B18").Value = "affiliate.marketing.guide"

# --- Row 19 (new) ---
$ws.Range("A19").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B19").Value = "com.sugar.powerfulquotes"

# --- Row heights (auto-computed by the original author's app from wrapped
# text, reproduced explicitly here). ---
$ws.Rows.Item(1).RowHeight = 12.8
$ws.Rows.Item(2).RowHeight = 46.5
$ws.Rows.Item(3).RowHeight = 24
$ws.Rows.Item(4).RowHeight = 12.8
$ws.Rows.Item(5).RowHeight = 24
$ws.Rows.Item(6).RowHeight = 46.5
$ws.Rows.Item(7).RowHeight = 12.8
$ws.Rows.Item(8).RowHeight = 24
$ws.Rows.Item(9).RowHeight = 24
$ws.Rows.Item(10).RowHeight = 46.5
$ws.Rows.Item(11).RowHeight = 12.8
$ws.Rows.Item(12).RowHeight = 12.8
$ws.Rows.Item(13).RowHeight = 12.8
$ws.Rows.Item(14).RowHeight = 24
$ws.Rows.Item(15).RowHeight = 12.8
$ws.Rows.Item(16).RowHeight = 12.8
$ws.Rows.Item(17).RowHeight = 12.8
$ws.Rows.Item(18).RowHeight = 24
$ws.Rows.Item(19).RowHeight = 46.5

# --- Selection / view state moves to the top of the sheet, cursor on B6. ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select() | Out-Null
